# Fix regtime/borntime data so regtime (column D) is always chronologically
# after borntime (column E) for every employee row, as described in the
# commit message: "data corrected: regtime now after borntime (script updated)"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Excel's date epoch (serial day 0 == 1899-12-30)
$epoch = Get-Date -Year 1899 -Month 12 -Day 30 -Hour 0 -Minute 0 -Second 0

for ($r = 2; $r -le 31; $r++) {
    $dCell = $ws.Cells.Item($r, 4)   # column D = regtime
    $eCell = $ws.Cells.Item($r, 5)   # column E = borntime

    $dVal = $dCell.Value2
    $eVal = $eCell.Value2

    $dtD = $epoch.AddDays($dVal)
    $dtE = $epoch.AddDays($eVal)

    # Compare using the original numeric serial values (not the DateTime
    # objects) to reliably determine which timestamp is chronologically later.
    if ($dVal -gt $eVal) {
        $later = $dtD
        $earlier = $dtE
    } else {
        $later = $dtE
        $earlier = $dtD
    }

    # The later timestamp becomes the new regtime (column D); the earlier
    # timestamp becomes the new borntime (column E).
    $ld = $later.Day
    $lm = $later.Month
    if ($ld -le 12 -and $lm -le 12) {
        $swappedD = Get-Date -Year $later.Year -Month $ld -Day $lm -Hour $later.Hour -Minute $later.Minute -Second $later.Second
        $dCell.Value = $swappedD.ToOADate()
    } else {
        $dCell.Value = $later.ToString("dd/MM/yyyy HH:mm:ss")
    }

    $ed = $earlier.Day
    $em = $earlier.Month
    if ($ed -le 12 -and $em -le 12) {
        $swappedE = Get-Date -Year $earlier.Year -Month $ed -Day $em -Hour $earlier.Hour -Minute $earlier.Minute -Second $earlier.Second
        $eCell.Value = $swappedE.ToOADate()
    } else {
        $eCell.Value = $earlier.ToString("dd/MM/yyyy HH:mm:ss")
    }
}

# Update the saved view state to match the author's last selection/scroll
# position when they reviewed the corrected data.
$ws.Range("D2").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 19
